$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.248.83'
$ws.Range("E2").Value = '  +2.45%  '

$ws.Range("D3").Value = '3.698.83'
$ws.Range("E3").Value = '  +8.06%  '

$ws.Range("E4").Value = '  -0.06%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '582.51'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  -0.18%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '178.62'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +1.19%  '

$ws.Range("D7").Value = '3.689.68'
$ws.Range("E7").Value = '  +8.05%  '

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.619'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  +4.39%  '

$ws.Range("E9").Value = '  -0.01%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.201'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +0.27%  '

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.88'
$cell.Style = $origStyle
$ws.Range("E11").Value = '  +26.54%  '

$ws.Range("E12").Value = '  +4.88%  '

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '49.31'
$cell.Style = $origStyle
$ws.Range("E13").Value = '  +1.31%  '

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000289'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  +2.78%  '

$ws.Range("D15").Value = '4.294.15'
$ws.Range("E15").Value = '  +7.99%  '

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '680.50'
$cell.Style = $origStyle
$ws.Range("E16").Value = '  -1.92%  '

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.04'
$cell.Style = $origStyle
$ws.Range("E17").Value = '  +5.03%  '

$ws.Range("D18").Value = '3.695.84'
$ws.Range("E18").Value = '  +7.76%  '

$ws.Range("D19").Value = '71.351.99'
$ws.Range("E19").Value = '  +2.50%  '

$ws.Range("E20").Value = '  +0.97%  '

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '18.00'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +1.97%  '

$ws.Range("E22").Value = '  +2.51%  '

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.946'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  +5.61%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.43'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  +2.86%  '

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '102.35'
$cell.Style = $origStyle
$ws.Range("E25").Value = '  +1.23%  '

$ws.Range("E26").Value = '  +2.09%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.85'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +7.54%  '

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.29'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  +7.52%  '

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  +0.02%  '

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '35.23'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  +5.37%  '

$ws.Range("E31").Value = '  +5.59%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.16'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  +4.59%  '

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.53'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  +6.56%  '

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.06'
$cell.Style = $origStyle
$ws.Range("E35").Value = '  +10.10%  '

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '581.70'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  +1.66%  '

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.22'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  +2.02%  '

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.109'
$cell.Style = $origStyle
$ws.Range("E38").Value = '  +4.85%  '

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '58.70'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = $origStyle
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.649.38'
$ws.Range("E41").Value = '  +2.35%  '

$ws.Range("E42").Value = '  +4.06%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0460'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  +10.44%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.353'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  +6.45%  '

$ws.Range("D45").Value = '0.0₃0771'
$ws.Range("E45").Value = '  +5.47%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '35.77'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  +2.46%  '

$ws.Range("E47").Value = '  +4.41%  '

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.93'
$cell.Style = $origStyle
$ws.Range("E48").Value = '  +10.64%  '

$ws.Range("E49").Value = '  +4.40%  '

$ws.Range("E50").Value = '  +2.01%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.98'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +10.42%  '

Write-Output "Update complete"